# Update NATMI LR-pair sheet (Ccl21b-Ccr7) with new TPM-derived numbers.
# The six sending/target-cluster combinations (FAPs/MuSCs x Inflammatory-Mac/
# Resolving-Mac/MuSCs) are rewritten in place; rows 6-7 are brand new.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A:T = Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then the 16 NATMI expression/specificity metric columns.
$rows = @(
  @("FAPs","Ccl21b","Ccr7","Inflammatory-Mac",3,1,0.3774076666666666,1.132223,
    0.4877525841056716,0.588184597482006,2,0.6666666666666666,2.005747666666667,
    6.017243,0.5060070458025943,0.5076214806746766,0.7569845467987777,
    6.812860921188999,0.2468062441658923,0.2985751362838546),

  @("FAPs","Ccl21b","Ccr7","MuSCs",3,1,0.3774076666666666,1.132223,
    0.4877525841056716,0.588184597482006,1,0.5,0.03782,0.07564,
    0.009541173493701741,0.006381076648929175,0.01427355795333333,
    0.08564134772,0.004653732026953562,0.003753251000252235),

  @("FAPs","Ccl21b","Ccr7","Resolving-Mac",3,1,0.3774076666666666,1.132223,
    0.4877525841056716,0.588184597482006,3,1,1.920305333333333,5.760916,
    0.484451780703704,0.4859974426763943,0.7247379551408888,6.522641596268,
    0.2362926079128257,0.2858562101978993),

  @("MuSCs","Ccl21b","Ccr7","Inflammatory-Mac",2,1,0.396361,0.7927219999999999,
    0.5122474158943284,0.411815402517994,2,0.6666666666666666,2.005747666666667,
    6.017243,0.5060070458025943,0.5076214806746766,0.7950001509076666,
    4.770000905446,0.259200801636702,0.209046344390822),

  @("MuSCs","Ccl21b","Ccr7","MuSCs",2,1,0.396361,0.7927219999999999,
    0.5122474158943284,0.411815402517994,1,0.5,0.03782,0.07564,
    0.009541173493701741,0.006381076648929175,0.01499037302,
    0.05996149207999999,0.004887441466748178,0.002627825648676941),

  @("MuSCs","Ccl21b","Ccr7","Resolving-Mac",2,1,0.396361,0.7927219999999999,
    0.5122474158943284,0.411815402517994,3,1,1.920305333333333,5.760916,
    0.484451780703704,0.4859974426763943,0.7611341422253333,4.566804853351999,
    0.2481591727908782,0.200141232478495)
)

# Leave the header row (row 1) untouched - it already matches.
$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $values[$c]
    }
}
